# Coke Cooler Purity scene/session hierarchy fully established and working
#
# The "Targets" tab used to key its scene-type rows off combined labels like
# "SOS SSD" / "SOS Still" / ... together with a single "UNITED" delivery
# program column that was folded into the scene-type text. This flattens the
# hierarchy: column A now holds the bare scene type (SSD, Still, Isotonic,
# ...), column B explicitly carries the "UNITED" program/session, and column
# C keeps the store-attribute (Baton Rouge Preferred/Partnership/Parity).
# The now-unused "SOS ..." labels are dropped from the shared string table
# automatically once nothing references them any more.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Targets")
    $ws.Range("A2").Value = "SSD"
    $ws.Range("B2").Value = "UNITED"
    $ws.Range("C2").Value = "Baton Rouge Preferred"
    $ws.Range("A3").Value = "Still"
    $ws.Range("B3").Value = "UNITED"
    $ws.Range("C3").Value = "Baton Rouge Preferred"
    $ws.Range("A4").Value = "Isotonic"
    $ws.Range("B4").Value = "UNITED"
    $ws.Range("C4").Value = "Baton Rouge Preferred"
    $ws.Range("A5").Value = "Tea"
    $ws.Range("B5").Value = "UNITED"
    $ws.Range("C5").Value = "Baton Rouge Preferred"
    $ws.Range("A6").Value = "Water"
    $ws.Range("B6").Value = "UNITED"
    $ws.Range("C6").Value = "Baton Rouge Preferred"
    $ws.Range("A7").Value = "Energy"
    $ws.Range("B7").Value = "UNITED"
    $ws.Range("C7").Value = "Baton Rouge Preferred"
    $ws.Range("A8").Value = "Juice"
    $ws.Range("B8").Value = "UNITED"
    $ws.Range("C8").Value = "Baton Rouge Preferred"
    $ws.Range("A9").Value = "Enhanced Water (Vitamin Water)"
    $ws.Range("B9").Value = "UNITED"
    $ws.Range("C9").Value = "Baton Rouge Preferred"
    $ws.Range("A10").Value = "SSD"
    $ws.Range("B10").Value = "UNITED"
    $ws.Range("C10").Value = "Baton Rouge Partnership"
    $ws.Range("A11").Value = "Still"
    $ws.Range("B11").Value = "UNITED"
    $ws.Range("C11").Value = "Baton Rouge Partnership"
    $ws.Range("A12").Value = "Isotonic"
    $ws.Range("B12").Value = "UNITED"
    $ws.Range("C12").Value = "Baton Rouge Partnership"
    $ws.Range("A13").Value = "Tea"
    $ws.Range("B13").Value = "UNITED"
    $ws.Range("C13").Value = "Baton Rouge Partnership"
    $ws.Range("A14").Value = "Water"
    $ws.Range("B14").Value = "UNITED"
    $ws.Range("C14").Value = "Baton Rouge Partnership"
    $ws.Range("A15").Value = "Energy"
    $ws.Range("B15").Value = "UNITED"
    $ws.Range("C15").Value = "Baton Rouge Partnership"
    $ws.Range("A16").Value = "Juice"
    $ws.Range("B16").Value = "UNITED"
    $ws.Range("C16").Value = "Baton Rouge Partnership"
    $ws.Range("A17").Value = "Enhanced Water (Vitamin Water)"
    $ws.Range("B17").Value = "UNITED"
    $ws.Range("C17").Value = "Baton Rouge Partnership"
    $ws.Range("A18").Value = "SSD"
    $ws.Range("B18").Value = "UNITED"
    $ws.Range("C18").Value = "Baton Rouge Parity"
    $ws.Range("A19").Value = "Still"
    $ws.Range("B19").Value = "UNITED"
    $ws.Range("C19").Value = "Baton Rouge Parity"
    $ws.Range("A20").Value = "Isotonic"
    $ws.Range("B20").Value = "UNITED"
    $ws.Range("C20").Value = "Baton Rouge Parity"
    $ws.Range("A21").Value = "Tea"
    $ws.Range("B21").Value = "UNITED"
    $ws.Range("C21").Value = "Baton Rouge Parity"
    $ws.Range("A22").Value = "Water"
    $ws.Range("B22").Value = "UNITED"
    $ws.Range("C22").Value = "Baton Rouge Parity"
    $ws.Range("A23").Value = "Energy"
    $ws.Range("B23").Value = "UNITED"
    $ws.Range("C23").Value = "Baton Rouge Parity"
    $ws.Range("A24").Value = "Juice"
    $ws.Range("B24").Value = "UNITED"
    $ws.Range("C24").Value = "Baton Rouge Parity"
    $ws.Range("A25").Value = "Enhanced Water (Vitamin Water)"
    $ws.Range("B25").Value = "UNITED"
    $ws.Range("C25").Value = "Baton Rouge Parity"

# --- cosmetic: column widths nudged slightly wider across the workbook ---
$wb.Worksheets.Item("KPIs").Columns.Item(1).ColumnWidth = 29.333333333333336
$wb.Worksheets.Item("KPIs").Columns.Item(3).ColumnWidth = 37.166666666666664
$wb.Worksheets.Item("KPIs").Columns.Item(4).ColumnWidth = 9.333333333333332
$wb.Worksheets.Item("Visible").Columns.Item(2).ColumnWidth = 47.666666666666664
$wb.Worksheets.Item("SOS").Columns.Item(1).ColumnWidth = 29.333333333333336
$wb.Worksheets.Item("SOS").Columns.Item(2).ColumnWidth = 17.333333333333336
$wb.Worksheets.Item("SOS").Columns.Item(3).ColumnWidth = 25.166666666666668
$wb.Worksheets.Item("SOS").Columns.Item(4).ColumnWidth = 17.333333333333336
$wb.Worksheets.Item("SOS").Columns.Item(6).ColumnWidth = 19.166666666666668
$wb.Worksheets.Item("SOS").Columns.Item(7).ColumnWidth = 18.166666666666668
$wb.Worksheets.Item("SOS").Columns.Item(8).ColumnWidth = 19.166666666666668
$wb.Worksheets.Item("SOS").Columns.Item(9).ColumnWidth = 18.166666666666668
$wb.Worksheets.Item("SOS").Columns.Item(13).ColumnWidth = 27.5
$wb.Worksheets.Item("Targets").Columns.Item(1).ColumnWidth = 29.333333333333336
$wb.Worksheets.Item("Targets").Columns.Item(3).ColumnWidth = 21.833333333333336

# The Targets sheet is now the one the user was last working in
$ws.Range("A26").Select()
$ws.Activate()
